$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string; force Text format first
# so Excel does not auto-convert/round them into actual numbers.
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '22.493.24'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.576.60'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '1.001'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = '288.71'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').Value = '0.3685'
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('D8').Value = '48.12'
$ws.Range('E8').Value = '  -3.25%  '
$ws.Range('D9').Value = '0.3331'
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('D10').Value = '1.149'
$ws.Range('E10').Value = '  +1.89%  '
$ws.Range('D11').Value = '0.07558'
$ws.Range('E11').Value = '  +2.07%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '20.84'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '5.973'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = '6.952'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').Value = '1.575.09'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D17').Value = '0.00001124'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').Value = '87.95'
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('D19').Value = '0.06734'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.399'
$ws.Range('E21').Value = '  +1.60%  '
$ws.Range('D22').Value = '16.62'
$ws.Range('E22').Value = '  +3.64%  '
$ws.Range('D23').Value = '12.01'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').Value = '22.492.82'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').Value = '2.390'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').Value = '2.647'
$ws.Range('E26').Value = '  +3.67%  '
$ws.Range('D27').Value = '151.34'
$ws.Range('E27').Value = '  +1.36%  '
$ws.Range('D28').Value = '19.70'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').Value = '5.006'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').Value = '125.43'
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('D31').Value = '1.751.51'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').Value = '1.102'
$ws.Range('E32').Value = '  +4.91%  '
$ws.Range('D33').Value = '6.155'
$ws.Range('E33').Value = '  +0.89%  '
$ws.Range('D34').Value = '1.996'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').Value = '9.890'
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('D36').Value = '0.08380'
$ws.Range('E36').Value = '  +1.28%  '
$ws.Range('D37').Value = '0.02477'
$ws.Range('E37').Value = '  +3.60%  '
$ws.Range('D38').Value = '0.2255'
$ws.Range('E38').Value = '  +1.70%  '
$ws.Range('D39').Value = '0.06424'
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('D40').Value = '5.384'
$ws.Range('E40').Value = '  +1.19%  '
$ws.Range('D41').Value = '1.295'
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('D42').Value = '11.49'
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('D43').Value = '0.6309'
$ws.Range('E43').Value = '  +3.75%  '
$ws.Range('D44').Value = '14.06'
$ws.Range('E44').Value = '  +2.33%  '
$ws.Range('D45').Value = '1.0000'
$ws.Range('D46').Value = '0.6145'
$ws.Range('E46').Value = '  +7.35%  '
$ws.Range('D47').Value = '3.791'
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('D48').Value = '2.075'
$ws.Range('E48').Value = '  +2.93%  '
$ws.Range('D49').Value = '125.84'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').Value = '1.217'
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('D51').Value = '0.07240'
$ws.Range('E51').Value = '  -0.03%  '
